{"js": "// The commit rewrites the last three bulleted \"Summary\" items:\n//   1. The \"Peer Review: ...\" bullet becomes the \"Search for Properties Test\n//      Cases: ...\" bullet (and the _GoBack bookmark moves into it, splitting\n//      \"of test cases...\" out into its own run).\n//   2. The (old) \"Search for Properties Test Cases: ...\" bullet becomes the\n//      \"register.php: ...\" bullet, now with proofErr spell-check wrappers\n//      around \"register.php\" and \"php\".\n//   3. The (old) \"register.php: ...\" bullet (which used to carry the\n//      _GoBack bookmark) becomes a brand-new \"addproperty.php: ...\" bullet,\n//      with proofErr wrappers, and no longer carries the bookmark.\n//\n// We rebuild each paragraph's contents precisely (runs, proofErr marks and\n// the bookmark) via Range.insertOoxml(..., \"Replace\") so the resulting\n// run/bookmark/proofErr structure matches the target exactly, rather than\n// relying on generic text replacement which can't create proofErr or move\n// bookmarks between paragraphs.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the three target paragraphs by their (pre-edit) text so the script\n// is resilient to any surrounding content.\nlet peerReviewIndex = -1;\nlet searchCasesIndex = -1;\nlet registerPhpIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Peer Review:\") === 0) peerReviewIndex = i;\n  if (t.indexOf(\"Search for Properties Test Cases\") === 0) searchCasesIndex = i;\n  if (t.indexOf(\"register.php\") === 0) registerPhpIndex = i;\n}\n\nif (peerReviewIndex === -1 || searchCasesIndex === -1 || registerPhpIndex === -1) {\n  throw new Error(\n    \"Could not locate expected paragraphs (peerReview=\" + peerReviewIndex +\n    \", searchCases=\" + searchCasesIndex + \", registerPhp=\" + registerPhpIndex + \")\"\n  );\n}\n\nfunction flatOpc(innerParagraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    innerParagraphXml +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst rPr = '<w:rPr><w:lang w:val=\"en-AU\"/></w:rPr>';\nconst pPr =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>' +\n  rPr +\n  '</w:pPr>';\n\n// New paragraph 1: \"Search for Properties Test Cases: A list [bookmark]of test cases...\"\nconst newPeerReviewParaXml =\n  '<w:p>' +\n  pPr +\n  '<w:r>' + rPr + '<w:t>Search for Properties Test Cases</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">: A list </w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r>' + rPr + '<w:t>of test cases that were carried out on one of our pages.  Similar test cases were done for each page</w:t></w:r>' +\n  '</w:p>';\n\n// New paragraph 2: \"register.php: This is an example of the php coding I have done in this unit.\"\nconst newSearchCasesParaXml =\n  '<w:p>' +\n  pPr +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>register.php</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">: This is an example of the </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>php</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\"> coding I have done in this unit.</w:t></w:r>' +\n  '</w:p>';\n\n// New paragraph 3: \"addproperty.php: This is an example of the php coding I have done in this unit.\" (no bookmark)\nconst newRegisterPhpParaXml =\n  '<w:p>' +\n  pPr +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>a</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t>ddproperty.php</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">This is an example of the </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>php</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\"> coding I have done in this unit.</w:t></w:r>' +\n  '</w:p>';\n\n// Apply replacements from bottom to top so earlier indices stay valid.\nconst ordered = [\n  { index: registerPhpIndex, xml: newRegisterPhpParaXml },\n  { index: searchCasesIndex, xml: newSearchCasesParaXml },\n  { index: peerReviewIndex, xml: newPeerReviewParaXml },\n].sort((a, b) => b.index - a.index);\n\nfor (const { index, xml } of ordered) {\n  const range = paragraphs.items[index].getRange(\"Whole\");\n  range.insertOoxml(flatOpc(xml), Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The commit rewrites the last three bulleted \"Summary\" items:\n#   1. The \"Peer Review: ...\" bullet becomes the \"Search for Properties Test\n#      Cases: ...\" bullet (and the _GoBack bookmark moves into it, splitting\n#      \"of test cases...\" out into its own run).\n#   2. The (old) \"Search for Properties Test Cases: ...\" bullet becomes the\n#      \"register.php: ...\" bullet, now with proofErr spell-check wrappers\n#      around \"register.php\" and \"php\".\n#   3. The (old) \"register.php: ...\" bullet (which used to carry the\n#      _GoBack bookmark) becomes a brand-new \"addproperty.php: ...\" bullet,\n#      with proofErr wrappers, and no longer carries the bookmark.\n#\n# Each paragraph is rebuilt precisely (runs, proofErr marks, bookmark) via\n# Range.InsertXML(...) -- the COM equivalent of Office.js's\n# Range.insertOoxml(..., \"Replace\") -- so the resulting run/bookmark/proofErr\n# structure matches the target exactly.\n\n$d = $word.ActiveDocument\n\nfunction New-FlatOpc([string]$innerParagraphXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        $innerParagraphXml +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n$rPr = '<w:rPr><w:lang w:val=\"en-AU\"/></w:rPr>'\n$pPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>' + $rPr + '</w:pPr>'\n\n# New paragraph 1: \"Search for Properties Test Cases: A list [bookmark]of test cases...\"\n$newPeerReviewParaXml = '<w:p>' + $pPr +\n    '<w:r>' + $rPr + '<w:t>Search for Properties Test Cases</w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">: A list </w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r>' + $rPr + '<w:t>of test cases that were carried out on one of our pages.  Similar test cases were done for each page</w:t></w:r>' +\n    '</w:p>'\n\n# New paragraph 2: \"register.php: This is an example of the php coding I have done in this unit.\"\n$newSearchCasesParaXml = '<w:p>' + $pPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' + $rPr + '<w:t>register.php</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">: This is an example of the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' + $rPr + '<w:t>php</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\"> coding I have done in this unit.</w:t></w:r>' +\n    '</w:p>'\n\n# New paragraph 3: \"addproperty.php: This is an example of the php coding I have done in this unit.\" (no bookmark)\n$newRegisterPhpParaXml = '<w:p>' + $pPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' + $rPr + '<w:t>a</w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t>ddproperty.php</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">: </w:t></w:r>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\">This is an example of the </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' + $rPr + '<w:t>php</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r>' + $rPr + '<w:t xml:space=\"preserve\"> coding I have done in this unit.</w:t></w:r>' +\n    '</w:p>'\n\n# Locate the three target paragraphs by their (pre-edit) text so the script\n# is resilient to any surrounding content. Record 1-based Paragraphs indices.\n$peerReviewIndex = -1\n$searchCasesIndex = -1\n$registerPhpIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Peer Review:\")) { $peerReviewIndex = $i }\n    if ($t.StartsWith(\"Search for Properties Test Cases\")) { $searchCasesIndex = $i }\n    if ($t.StartsWith(\"register.php\")) { $registerPhpIndex = $i }\n}\n\nif ($peerReviewIndex -eq -1 -or $searchCasesIndex -eq -1 -or $registerPhpIndex -eq -1) {\n    throw \"Could not locate expected paragraphs (peerReview=$peerReviewIndex, searchCases=$searchCasesIndex, registerPhp=$registerPhpIndex)\"\n}\n\n# Apply replacements from the highest paragraph index down to the lowest so\n# earlier indices stay valid as we go.\n$ordered = @(\n    @{ Index = $registerPhpIndex; Xml = $newRegisterPhpParaXml },\n    @{ Index = $searchCasesIndex; Xml = $newSearchCasesParaXml },\n    @{ Index = $peerReviewIndex; Xml = $newPeerReviewParaXml }\n) | Sort-Object -Property Index -Descending\n\nforeach ($item in $ordered) {\n    $range = $d.Paragraphs.Item($item.Index).Range\n    $range.InsertXML((New-FlatOpc $item.Xml))\n}\n"}
